$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText -replace "6605\.65 pesos", "6597.13 pesos" -replace "968\.41 Bs", "956.87 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 551.3
$wsTasas.Range("O10").Value = 3637
$wsTasas.Range("N12").Value = 3654.07
$wsTasas.Range("O12").Value = 530
